# "Generate Report for Handoff"
# - Status moves from "In Translation" to "Ready for handoff" for both locales.
# - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#   advance a couple of minutes to reflect the new handoff generation run.
# - The Status column is widened (in Overview + each locale sheet) to fit the
#   new, longer "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Handoff generation timestamps ---------------------------------------
$overview.Range("G2").Value = "2016-12-15 05:02:02"
$zhcn.Range("H2").Value = "2016-12-15 05:01:49"
$dede.Range("H2").Value = "2016-12-15 05:02:02"

# --- Widen the status columns to fit "Ready for handoff" -----------------
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3
